# Adding the changes we made on may 9th
# - Insert 4 new accelerometer samples at the top of the data (right after the
#   header row), pushing the existing rows down.
# - Append 6 new accelerometer samples at the end of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: make room for the 4 new rows. Inserting at rows 3:6 (instead of
# 2:5) means the newly created blank rows inherit the plain/default
# formatting of row 2 rather than the bold/bordered header formatting of
# row 1, which keeps the data rows styled the same way as before.
$ws.Range("A3:C6").EntireRow.Insert()

# Step 2: the original row 2 values are still sitting in row 2 (rows 3-6
# are the new blank ones) - move them down to row 6 so rows 2-5 are free
# for the new data, restoring the original row order (old row2 -> row6,
# old row3 -> row7, ... which already happened automatically for rows 7+).
$ws.Range("A2:C2").Cut($ws.Range("A6:C6"))

# Step 3: fill in the 4 new rows at the top (new rows 2-5).
$ws.Range("A2").Value = -1.350502490997314
$ws.Range("B2").Value = -5.659398555755615
$ws.Range("C2").Value = -3.780858278274536

$ws.Range("A3").Value = 6.421082973480225
$ws.Range("B3").Value = -10.34209537506104
$ws.Range("C3").Value = -8.60602855682373

$ws.Range("A4").Value = -8.985816955566406
$ws.Range("B4").Value = -4.653406143188477
$ws.Range("C4").Value = -4.82445764541626

$ws.Range("A5").Value = 8.638599395751953
$ws.Range("B5").Value = -36.5815658569336
$ws.Range("C5").Value = 4.202666282653809

# Step 4: append the 6 new rows at the bottom (rows 26-31).
$ws.Range("A26").Value = -12.95152854919434
$ws.Range("B26").Value = -7.215466022491455
$ws.Range("C26").Value = 2.712479591369629

$ws.Range("A27").Value = -1.57819402217865
$ws.Range("B27").Value = 4.682579040527344
$ws.Range("C27").Value = -4.256365299224854

$ws.Range("A28").Value = -8.448655128479004
$ws.Range("B28").Value = -33.68074035644531
$ws.Range("C28").Value = -12.51845359802246

$ws.Range("A29").Value = 71.44110107421875
$ws.Range("B29").Value = 6.458947658538818
$ws.Range("C29").Value = -16.25710296630859

$ws.Range("A30").Value = 7.958407878875732
$ws.Range("B30").Value = -4.263494491577148
$ws.Range("C30").Value = -17.78672790527344

$ws.Range("A31").Value = -30.355712890625
$ws.Range("B31").Value = -23.57818984985352
$ws.Range("C31").Value = 8.899335861206055
